# Dev #268: spec for new tab model
# Rename the single worksheet from "census" to "Tab". Excel automatically
# rewrites any defined names / formulas that reference the sheet by name
# (e.g. the hidden _FilterDatabase name used by the autofilter), so no
# extra work is needed there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Tab"

# Reset the saved selection on the sheet to a plain A2 (it previously
# covered the whole sheet, "1:1048576").
$ws.Range("A2").Select()

# Nudge column J's width very slightly (13.28 -> 13.29 characters).
$ws.Columns.Item(10).ColumnWidth = 13.29
